# Automatische test-sync: 2025-06-19 22:37:50
# Adds a new mail log entry to the "Logs" sheet and updates the
# "Dashboard" category-count table (and associated conditional
# formatting ranges) to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append the new mail entry as row 50
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(50, 1).Value = "Offerte voor zakelijke samenwerking"
$logs.Cells.Item(50, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(50, 3).Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Cells.Item(50, 4).Value = "Offerte / Prijsaanvraag"
$logs.Cells.Item(50, 6).Value = "2025-06-19 22:37:14"
$logs.Cells.Item(50, 7).Value = "Nee"

# Extend the conditional formatting ranges (D2:D49 -> D2:D50 and
# G2:G49 -> G2:G50) so the new row is covered, while preserving the
# existing rules/colors exactly.
$dRange = $logs.Range("D2:D49")
$dRange.FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D50"))

$gRange = $logs.Range("G2:G49")
$gRange.FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G50"))

# ---------------------------------------------------------------
# 2. Dashboard sheet: the new entry's category "Offerte /
#    Prijsaanvraag" now has 6 occurrences, overtaking "Afmelding /
#    Nieuwsbrief" (5) in the count-sorted table, so rows 4 and 5
#    swap places.
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(4, 2).Value = 6
$dash.Cells.Item(5, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(5, 2).Value = 5
